$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume(1h) figures. Numeric-looking D-column
# values are assigned with a leading apostrophe so Excel keeps them as text
# (matching the sheet's existing inlineStr cell type) instead of coercing
# them to numbers.
$ws.Range("D2").Value = '67.288.49'
$ws.Range("E2").Value = '  +4.75%  '
$ws.Range("D3").Value = '3.483.22'
$ws.Range("E3").Value = '  +4.81%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''584.96'
$ws.Range("E5").Value = '  +5.84%  '
$ws.Range("D6").Value = '''185.72'
$ws.Range("E6").Value = '  +7.36%  '
$ws.Range("D7").Value = '''0.634'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").Value = '3.477.33'
$ws.Range("E8").Value = '  +4.87%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '''0.173'
$ws.Range("E10").Value = '  +1.66%  '
$ws.Range("D11").Value = '''0.651'
$ws.Range("E11").Value = '  +2.94%  '
$ws.Range("D12").Value = '''56.27'
$ws.Range("E12").Value = '  +5.96%  '
$ws.Range("E13").Value = '  +1.05%  '
$ws.Range("D14").Value = '''9.46'
$ws.Range("E14").Value = '  +4.38%  '
$ws.Range("D15").Value = '4.045.79'
$ws.Range("E15").Value = '  +4.85%  '
$ws.Range("D16").Value = '''18.86'
$ws.Range("E16").Value = '  +4.25%  '
$ws.Range("D17").Value = '3.492.41'
$ws.Range("E17").Value = '  +4.89%  '
$ws.Range("D18").Value = '67.399.63'
$ws.Range("E18").Value = '  +4.72%  '
$ws.Range("D19").Value = '''12.19'
$ws.Range("E19").Value = '  +4.03%  '
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  +3.72%  '
$ws.Range("D22").Value = '''489.34'
$ws.Range("E22").Value = '  +8.21%  '
$ws.Range("D23").Value = '''5.39'
$ws.Range("E23").Value = '  +9.25%  '
$ws.Range("D24").Value = '''16.86'
$ws.Range("E24").Value = '  +21.63%  '
$ws.Range("D25").Value = '''4.44'
$ws.Range("E25").Value = '  +9.55%  '
$ws.Range("D26").Value = '''89.92'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("D27").Value = '''2.95'
$ws.Range("E27").Value = '  +2.80%  '
$ws.Range("D28").Value = '''11.01'
$ws.Range("E28").Value = '  +4.55%  '
$ws.Range("D29").Value = '''9.16'
$ws.Range("E29").Value = '  +6.86%  '
$ws.Range("D30").Value = '''31.49'
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("D31").Value = '''7.23'
$ws.Range("E31").Value = '  +11.00%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '''599.82'
$ws.Range("E32").Value = '  +5.02%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '''11.78'
$ws.Range("E33").Value = '  +3.40%  '
$ws.Range("D34").Value = '''63.87'
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("E35").Value = '  +4.96%  '
$ws.Range("E36").Value = '  +6.33%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '''36.62'
$ws.Range("E38").Value = '  +3.74%  '
$ws.Range("E39").Value = '  +0.97%  '
$ws.Range("D40").Value = '''0.387'
$ws.Range("E40").Value = '  +5.64%  '
$ws.Range("D41").Value = '0.0₃0764'
$ws.Range("E41").Value = '  +5.19%  '
$ws.Range("D42").Value = '3.269.98'
$ws.Range("E42").Value = '  +6.84%  '
$ws.Range("D43").Value = '''2.92'
$ws.Range("E43").Value = '  +6.87%  '
$ws.Range("E44").Value = '  +4.28%  '
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("E46").Value = '  +23.99%  '
$ws.Range("D47").Value = '''3.26'
$ws.Range("E47").Value = '  +2.27%  '
$ws.Range("E48").Value = '  +1.27%  '
$ws.Range("E49").Value = '  +12.27%  '
$ws.Range("D50").Value = '''8.77'
$ws.Range("E50").Value = '  +7.99%  '
$ws.Range("E51").Value = '  -0.04%  '
